$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 (the "Ponte Preta vs Sport Recife" Brazil Serie B game);
# this shifts row 3 (Colorado Springs vs Las Vegas Lights) up to become the new row 2.
$ws.Rows.Item(2).Delete()
